$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# Fill in the previously-empty "nan" placeholder cells on row 17 (B..K and P..R)
$ws.Range("B17").Value = "nan"
$ws.Range("C17").Value = "nan"
$ws.Range("D17").Value = "nan"
$ws.Range("E17").Value = "nan"
$ws.Range("F17").Value = "nan"
$ws.Range("G17").Value = "nan"
$ws.Range("H17").Value = "nan"
$ws.Range("I17").Value = "nan"
$ws.Range("J17").Value = "nan"
$ws.Range("K17").Value = "nan"
$ws.Range("P17").Value = "nan"
$ws.Range("Q17").Value = "nan"
$ws.Range("R17").Value = "nan"

# Add the new event row (row 18) for Card15
$ws.Range("A18").Value = "15"
$ws.Range("L18").Value = "12\8\2025"
$ws.Range("N18").Value = "تم تغييرزيت  الجيربوكس وتغيير جريد 1"
$ws.Range("O18").Value = "تيم العمل"
